$d = $word.ActiveDocument

$replacements = @(
    @("2026-01-31 Saturday", "2026-02-01 Sunday"),
    @("78×11=858", "65×77=5005"),
    @("11×36=396", "85×96=8160"),
    @("69×97=6693", "78×71=5538"),
    @("23×49=1127", "51×53=2703"),
    @("33×34=1122", "13×75=975"),
    @("87×30=2610", "11×15=165"),
    @("18×45=810", "13×43=559"),
    @("46×40=1840", "97×47=4559"),
    @("26×11=286", "94×85=7990"),
    @("68×79=5372", "84×97=8148"),
    @("70×80=5600", "16×90=1440"),
    @("64×13=832", "79×39=3081"),
    @("34×25=850", "39×58=2262"),
    @("66×17=1122", "94×36=3384"),
    @("53×19=1007", "69×69=4761"),
    @("36×64=2304", "50×39=1950"),
    @("65×53=3445", "31×25=775"),
    @("71×38=2698", "70×97=6790"),
    @("84×30=2520", "80×67=5360"),
    @("25×46=1150", "94×21=1974"),
    @("18×43=774", "67×51=3417"),
    @("32×50=1600", "80×21=1680"),
    @("88×89=7832", "33×38=1254"),
    @("19×74=1406", "26×26=676"),
    @("19×75=1425", "23×65=1495")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
